# Apply the "1% -> 5%" modifier change for units-vs-units on the UnitMod sheet.
# The workbook has a single worksheet ("UnitMod") which contains a handful of
# hard-coded 0.01 values (R/S columns for rows 21, 22, 26, 29). All of the
# other cells touched by the diff are formulas that derive from these base
# cells (directly or indirectly), so updating just the base cells and letting
# Excel recalculate reproduces the rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UnitMod")

# Update the base "units vs units" 1% modifiers to 5% (R/S columns).
$ws.Range("R21").Value = 0.05
$ws.Range("S21").Value = 0.05

$ws.Range("R22").Value = 0.05
$ws.Range("S22").Value = 0.05

$ws.Range("R26").Value = 0.05
$ws.Range("S26").Value = 0.05

$ws.Range("R29").Value = 0.05
$ws.Range("S29").Value = 0.05

# Force a full recalculation so all the dependent formula cells
# (rows 43, 44, 48, 62, 63, 67, 80, 81, 85, 88, 99, 100, 104, 105,
# 116, 121, 133, 137, 140, ...) pick up the new values.
$excel.CalculateFullRebuild()

# Update the selection to match the saved view state in the diff.
$ws.Range("R29").Select()

$wb.Save()
